$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 'This paper introduces SQUID (Surrogate Quantitative Interpretability for Deepnets), a framework for interpreting genomic deep neural networks by approximating them locally with simpler, inherently interpretable surrogate models that incorporate domain-specific knowledge about cis-regulatory mechanisms. By accounting for the confounding effects of nonlinearities and heteroscedastic noise in functional genomics data, SQUID outperforms existing interpretability methods in identifying consistent motifs across genomic loci and predicting variant effects. The framework also enables quantification of epistatic interactions within and between regulatory elements and provides global mechanistic explanations across sequence contexts, advancing the mechanistic interpretation of genomic DNNs.'
$ws.Range("G6").Value = 'Here we introduce MAVE-NN, a neural-network-based Python package that implements a broadly applicable information-theoretic framework for learning genotype-phenotype maps—including biophysically interpretable models—from MAVE datasets. We demonstrate MAVE-NN in multiple biological contexts, and highlight the ability of our approach to deconvolve mutational effects from otherwise confounding experimental nonlinearities and noise.'
$ws.Range("G8").Value = 'This review examines how massively parallel assays—including deep mutational scanning, high-throughput SELEX, and massively parallel reporter assays—have transformed the quantitative modeling of sequence–function relationships across diverse biological contexts, from clinical variant interpretation to transcription factor binding, protein landscapes, and cis-regulatory mechanisms. We present a unified conceptual framework and core mathematical modeling strategies applicable across these areas, spanning topics such as protein evolution, transcriptional regulation, and mRNA splicing. We emphasize critical principles of experimental design and mathematical modeling necessary for ensuring interpretability and reproducibility in such studies.'
$ws.Range("G9").Value = 'Here we describe a field-theoretic approach that addresses this problem remarkably well in one dimension, providing an exact nonparametric Bayesian posterior without relying on tunable parameters or large-data approximations. Strong non-Gaussian constraints, which require a nonperturbative treatment, are found to play a major role in reducing distribution uncertainty. A software implementation of this method is provided.'
$ws.Range("G10").Value = 'Here we describe a new experimental approach, called Tite-Seq, that is capable of measuring binding titration curves and corresponding affinities for thousands of variant antibodies in parallel. The measurement of titration curves eliminates the confounding effects of antibody expression and stability that arise in standard deep mutational scanning assays.'
$ws.Range("G11").Value = 'This paper proposes a mathematical formalization of "equitability"—the ability to quantify statistical associations without bias toward specific relationship forms—using core concepts from information theory. We show that mutual information, a fundamental information-theoretic measure of dependence, naturally satisfies this equitability criterion, whereas the recently introduced maximal information coefficient violates it. We conclude that estimating mutual information provides a natural and practical method for equitably quantifying associations in large datasets.'
$ws.Range("F12").Value = '/pictures/figures/fig1_deeplearning_kinney2010.jpeg'

$ws.Activate()
$ws.Range("F11").Select()
